# Update cryptocurrency price/volume data (scheduled GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.866.69"
$ws.Range("E2").Value = "  +0.31%  "

$ws.Range("D3").Value = "1.706.88"
$ws.Range("E3").Value = "  +0.32%  "

$ws.Range("D4").Value = "'0.9993"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.52%  "

$ws.Range("D5").Value = "'317.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.13%  "

$ws.Range("E6").Value = "  -0.57%  "

$ws.Range("D7").Value = "'0.3955"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "'0.4061"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.53%  "

$ws.Range("D9").Value = "'1.486"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.29%  "

$ws.Range("D10").Value = "'0.9994"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.54%  "

$ws.Range("D11").Value = "'53.22"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.67%  "

$ws.Range("D12").Value = "'0.08819"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.06%  "

$ws.Range("D13").Value = "'26.49"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +9.24%  "

$ws.Range("E14").Value = "  -2.56%  "

$ws.Range("D15").Value = "'8.138"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.08%  "

$ws.Range("D16").Value = "'0.00001361"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.19%  "

$ws.Range("D17").Value = "1.751.98"
$ws.Range("E17").Value = "  +2.48%  "

$ws.Range("D18").Value = "'96.63"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.20%  "

$ws.Range("D19").Value = "'0.07157"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.21%  "

$ws.Range("D20").Value = "'21.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.65%  "

$ws.Range("D21").Value = "'7.293"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.08%  "

$ws.Range("D22").Value = "'0.9986"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.64%  "

$ws.Range("D23").Value = "'14.33"
$ws.Range("D23").Style = "Normal"

$ws.Range("D24").Value = "24.827.88"
$ws.Range("E24").Value = "  +0.19%  "

$ws.Range("D25").Value = "'2.990"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.43%  "

$ws.Range("D26").Value = "'2.338"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.05%  "

$ws.Range("D27").Value = "'23.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.87%  "

$ws.Range("D28").Value = "'6.240"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +20.07%  "

$ws.Range("D29").Value = "'166.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.04%  "

$ws.Range("D30").Value = "'145.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.20%  "

$ws.Range("D31").Value = "'8.470"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -9.31%  "

$ws.Range("D32").Value = "1.926.58"
$ws.Range("E32").Value = "  +1.50%  "

$ws.Range("D33").Value = "'2.255"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +14.68%  "

$ws.Range("D34").Value = "'0.08801"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.16%  "

$ws.Range("D35").Value = "'0.03219"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.62%  "

$ws.Range("E36").Value = "  -10.83%  "

$ws.Range("D37").Value = "'1.032"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.52%  "

$ws.Range("D38").Value = "'0.2873"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.17%  "

$ws.Range("D39").Value = "'10.94"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.14%  "

$ws.Range("E40").Value = "  +6.91%  "

$ws.Range("D41").Value = "'0.09261"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.49%  "

$ws.Range("D42").Value = "'14.10"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.13%  "

$ws.Range("D43").Value = "'1.477"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.00%  "

$ws.Range("D44").Value = "'17.43"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.69%  "

$ws.Range("D45").Value = "'2.692"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.17%  "

$ws.Range("D46").Value = "'0.7403"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.89%  "

$ws.Range("D47").Value = "'4.242"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.09%  "

$ws.Range("D48").Value = "'1.392"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.65%  "

$ws.Range("D49").Value = "'0.9990"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.47%  "

$ws.Range("D50").Value = "'141.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.11%  "

$ws.Range("D51").Value = "'0.08346"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.43%  "
